$wb = $excel.ActiveWorkbook

$metaSheet = $wb.Worksheets.Item("SwateTemplateMetadata")
$metaSheet.Name = "isa_template"
$metaSheet.Activate()

$wb.Styles.Item("Hyperlink").Name = "Link"
$wb.Styles.Item("Normal").Name = "Standard"
